# Update localization status report: files "02be0b45-47b6-48c6-b03e-2682590bef7f.md"
# and "16a95424-8127-4d59-b9a6-c1cc32567089.md" moved from "Ready for handoff"
# to "In Translation" for both locales (zh-cn, de-de). Reflect the change on the
# Overview sheet as well as the per-locale detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "In Translation"
$overview.Range("C3").Value = "In Translation"
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
